$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell U1: "2025-11-22" as literal text (not auto-converted to a date),
# styled to match the other header cells (bold, bordered, centered).
$ws.Range("U1").Formula = "=""2025-11-22"""
$ws.Range("U1").Copy()
$ws.Range("U1").PasteSpecial(-4163)   # xlPasteValues: freeze formula result as a literal value
$ws.Range("T1").Copy()
$ws.Range("U1").PasteSpecial(-4122)   # xlPasteFormats: copy the header style from T1

# Rows 2-9: mark the new date column with the same "absent" marker used elsewhere,
# and bump the Total (S) column by 1 to account for the extra tracked day.
$ws.Range("U2").Value = "❌"
$ws.Range("S2").Value = 16

$ws.Range("U3").Value = "❌"
$ws.Range("S3").Value = 16

$ws.Range("U4").Value = "❌"
$ws.Range("S4").Value = 16

$ws.Range("U5").Value = "❌"
$ws.Range("S5").Value = 16

$ws.Range("U6").Value = "❌"
$ws.Range("S6").Value = 16

$ws.Range("U7").Value = "❌"
$ws.Range("S7").Value = 16

$ws.Range("U8").Value = "❌"
$ws.Range("S8").Value = 1

$ws.Range("U9").Value = "❌"
$ws.Range("S9").Value = 16
